$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- E8: replace shared string "Sanity check" with numeric value ---
$ws.Range("E8").Value = 0.23200000000000001

# --- A9: label shifts from "Two above cvombined" stays the same text (shared string index shifts after "Sanity check" removal) ---
$ws.Range("A9").Value = "Two above cvombined"

# --- E9: fix formula so it no longer errors ---
$ws.Range("E9").Formula = "=E7+E8-E6"

# --- New cell Q4 ---
$ws.Range("Q4").Formula = "=60*K4"

# --- New cell S9 ---
$ws.Range("S9").Formula = "=87/126*3/2"

# --- New cell S11 ---
$ws.Range("S11").Formula = "=73/87*4/3"

# --- New cell S13 ---
$ws.Range("S13").Formula = "=3/4"

# --- New cell S14 ---
$ws.Range("S14").Formula = "=73/87"

# --- Update selection to match target state ---
$ws.Range("N11").Select()

$wb.Save()
